# LOG.xlsx update — add three new logbook entries (rows 44-46) covering
# work on classes/unit tests, the mid-term presentation, and the failed
# pyqt5-in-vscode attempt / test report templates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date/start/end-time formatting from the last filled row (43)
# down onto the three newly-used rows so the new cells pick up the same
# styles (date format, time format) instead of the blank "s=1" style.
$ws.Range("B43:D43").Copy() | Out-Null
$ws.Range("B44:D46").PasteSpecial(-4122) | Out-Null

# Same for the "Bezigheden" (activity) column formatting.
$ws.Range("E43").Copy() | Out-Null
$ws.Range("E44:E46").PasteSpecial(-4122) | Out-Null

$ws.Application.CutCopyMode = $false

# Row 44 — 2022-11-02, 08:45 - 17:00
$ws.Range("B44").Value = 44867
$ws.Range("C44").Value = 0.36458333333333331
$ws.Range("D44").Value = 0.70833333333333337
$ws.Range("E44").Value = "Created a lot of classes, rewrote classes, wrote unit tests, ended with steppermotor class. "

# Row 45 — 2022-11-03, 09:00 - 15:00
$ws.Range("B45").Value = 44868
$ws.Range("C45").Value = 0.375
$ws.Range("D45").Value = 0.625
$ws.Range("E45").Value = "Worked more on software, fixed motor class unit tests and held my mid term presentation"

# Row 46 — 2022-11-04, 11:00 - 15:00
$ws.Range("B46").Value = 44869
$ws.Range("C46").Value = 0.45833333333333331
$ws.Range("D46").Value = 0.625
$ws.Range("E46").Value = "tried to get pyqt5 working in vsc (failed), finished creating test report templates"

# Move the active selection/view to the newly added last row, matching
# where the author was working when the file was saved.
$ws.Range("E46").Select() | Out-Null
